# Update "想去人数" (column F) figures across the four sheets to reflect
# freshly re-scraped counts (gh-pages data refresh at commit 456a3b4).

function Set-FValues($SheetName, $RowToValue) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowToValue.Keys) {
        $ws.Range("F$row").Value = $RowToValue[$row]
    }
}

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
Set-FValues "展览" @{
    2  = 1733
    3  = 10111
    6  = 594
    8  = 1609
    9  = 170
    10 = 376
    12 = 201
    15 = 1174
    17 = 31
    19 = 88
    20 = 345
    25 = 695
    31 = 398
    34 = 526
    35 = 602
    36 = 736
}

# 演出 (Performances)
Set-FValues "演出" @{
    10 = 6
    11 = 57
    20 = 603
    22 = 322
    23 = 683
    28 = 363
    31 = 207
    33 = 27
    35 = 186
    39 = 130
    40 = 5
    41 = 66
}

# 本地生活 (Local Life)
Set-FValues "本地生活" @{
    5  = 191
    6  = 2512
    7  = 4059
    10 = 290
}

# 全部类型 (All Types)
Set-FValues "全部类型" @{
    2  = 1733
    4  = 10111
    5  = 191
    7  = 4059
    9  = 290
    10 = 290
    11 = 594
    12 = 1609
    13 = 170
    14 = 376
    15 = 201
    17 = 6
    18 = 1174
    20 = 57
    24 = 88
    26 = 345
    29 = 322
    31 = 695
    34 = 363
    35 = 398
    38 = 526
    39 = 602
    40 = 207
    41 = 736
    47 = 130
    50 = 66
}
